$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1047, shifting the existing rows 1047:1148 down to 1049:1150
$ws.Rows.Item(1047).Resize(2).Insert()

# Fill in the data for the two newly inserted rows (1047 and 1048)
$ws.Range("A1047").Value = 4
$ws.Range("B1047").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C1047").Value = "Los Lagos"
$ws.Range("D1047").Value = 45132
$ws.Range("E1047").Value = 10
$ws.Range("F1047").Value = 100112020
$ws.Range("G1047").Value = "Tomate"
$ws.Range("H1047").Value = "Larga vida"
$ws.Range("I1047").Value = "Primera"
$ws.Range("J1047").Value = 300
$ws.Range("K1047").Value = 30000
$ws.Range("L1047").Value = 30000
$ws.Range("M1047").Value = 30000
$ws.Range("N1047").Value = "$/bandeja 20 kilos"
$ws.Range("O1047").Value = "Región de Arica y Parinacota"
$ws.Range("P1047").Value = 1500
$ws.Range("Q1047").Value = 20
$ws.Range("R1047").Value = "Hortaliza"

$ws.Range("A1048").Value = 4
$ws.Range("B1048").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C1048").Value = "Los Lagos"
$ws.Range("D1048").Value = 45132
$ws.Range("E1048").Value = 10
$ws.Range("F1048").Value = 100112020
$ws.Range("G1048").Value = "Tomate"
$ws.Range("H1048").Value = "Larga vida"
$ws.Range("I1048").Value = "Segunda"
$ws.Range("J1048").Value = 250
$ws.Range("K1048").Value = 26000
$ws.Range("L1048").Value = 26000
$ws.Range("M1048").Value = 26000
$ws.Range("N1048").Value = "$/bandeja 20 kilos"
$ws.Range("O1048").Value = "Región de Arica y Parinacota"
$ws.Range("P1048").Value = 1300
$ws.Range("Q1048").Value = 20
$ws.Range("R1048").Value = "Hortaliza"
